{"js": "// Apply the Zulu proofreading fixes described in the commit diff.\n// Each fix is a small, unique, in-run text substitution, so we use\n// Body.search() to locate the exact phrase and Range.insertText(...,\n// Word.InsertLocation.replace) to swap it in place without disturbing\n// run formatting.\nconst body = context.document.body;\n\nconst replacements = [\n  {\n    find: \"sicela uthumele i-email ithimbeni locwaningo ku-\",\n    replace: \"sicela uthumele i-email ithimba locwaningo ku-\"\n  },\n  {\n    find: \"kuqinisekise ukuthi uzizwa ukhululekile uma uphendula imibuzo.\",\n    replace: \"kuqinisekise ukuthi uzizwe ukhululekile uma uphendula imibuzo.\"\n  },\n  {\n    find: \"kanti iMenenja yocwaningo nguZamakhanya Makhanya (University of Cape Town).\",\n    replace: \"kanye neMenenja yocwaningo uZamakhanya Makhanya (University of Cape Town).\"\n  },\n  {\n    find: \"Uma unemibuzo noma okukukhathazayo mayelana namalungelo akho njengomhlanganyeli wocwaningo, ungathintana nethimba locwaningo ku-\",\n    replace: \"Uma unemibuzo noma kukhona okukukhathazayo mayelana namalungelo akho njengomhlanganyeli wocwaningo, ungathintana nethimba locwaningo ku-\"\n  },\n  {\n    find: \"Uma ufunde futhi waqonda idokhumenti engenhla, vuma kulemilayezo\",\n    replace: \"Uma ufunde futhi waqonda incwadi engenhla, vuma kulemilayezo\"\n  }\n];\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + find);\n  }\n\n  results.items[0].insertText(replace, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Apply the Zulu proofreading fixes described in the commit diff.\n# Each fix is a small, unique text substitution; we drive Word's\n# Find/Replace engine (Range.Find) over the whole document story for\n# each one so formatting/runs outside the matched text are left alone.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $found = $range.Find.Execute(\n        $findText,    # FindText\n        $true,        # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        1,            # Wrap (wdFindContinue)\n        $false,       # Format\n        $replaceText, # ReplaceWith\n        2             # Replace (wdReplaceOne)\n    )\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\nReplace-Text \"sicela uthumele i-email ithimbeni locwaningo ku-\" \"sicela uthumele i-email ithimba locwaningo ku-\"\nReplace-Text \"kuqinisekise ukuthi uzizwa ukhululekile uma uphendula imibuzo.\" \"kuqinisekise ukuthi uzizwe ukhululekile uma uphendula imibuzo.\"\nReplace-Text \"kanti iMenenja yocwaningo nguZamakhanya Makhanya (University of Cape Town).\" \"kanye neMenenja yocwaningo uZamakhanya Makhanya (University of Cape Town).\"\nReplace-Text \"Uma unemibuzo noma okukukhathazayo mayelana namalungelo akho njengomhlanganyeli wocwaningo, ungathintana nethimba locwaningo ku-\" \"Uma unemibuzo noma kukhona okukukhathazayo mayelana namalungelo akho njengomhlanganyeli wocwaningo, ungathintana nethimba locwaningo ku-\"\nReplace-Text \"Uma ufunde futhi waqonda idokhumenti engenhla, vuma kulemilayezo\" \"Uma ufunde futhi waqonda incwadi engenhla, vuma kulemilayezo\"\n"}
